$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that changed from 45205 (2023-10-06)
# to 45206 (2023-10-07) for every data row (rows 2 through 97).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
